# Update master data with new DAF values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artikel FGR+")

# Update DAF (column C) values for the rows that changed
$ws.Range("C16").Value = 1.2
$ws.Range("C18").Value = 1.2
$ws.Range("C21").Value = 1.9
$ws.Range("C23").Value = 2.5
$ws.Range("C24").Value = 1.4
$ws.Range("C25").Value = 1.8
$ws.Range("C27").Value = 1.7
$ws.Range("C28").Value = 1.6
$ws.Range("C29").Value = 1.8
$ws.Range("C30").Value = 2.6
$ws.Range("C31").Value = 1.4
$ws.Range("C33").Value = 2.3

# Update the view: scroll so row 13 is the top visible row, and move the
# active selection to C36
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("C36").Select()
$excel.ActiveWindow.ActivePane.ScrollRow = 13
